$d = $word.ActiveDocument
$rng = $d.Content

$xml = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p/><w:p/><w:p/><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>roundNumber</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}}</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>roundDate</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}}</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>{#</w:t></w:r><w:r><w:t>m</w:t></w:r><w:r><w:t>atches}</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>Zápas {{index}}</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>match_time</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}}</w:t></w:r><w:r><w:t xml:space="preserve"> -</w:t></w:r><w:r><w:t xml:space="preserve"> {{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>teamA</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}} x {{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>teamB</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}}</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>{/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>m</w:t></w:r><w:r><w:t>atches</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$rng.InsertXML($xml)
